$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "DFbVx581"
$ws.Range("B2").Value = 23091209
$ws.Range("C2").Value = "dzhujzq74"
$ws.Range("D2").Value = "Yk9&4jZ#"
$ws.Range("F2").Value = "YGScBgJI"
$ws.Range("G2").Value = "TJVJ"
